# Add "name" column to TextFileSequence sheet
$wb = $excel.ActiveWorkbook

$wsTextFileSequence = $wb.Worksheets.Item("TextFileSequence")
$wsTextFileSequence.Range("G1").Value = "name"

# Reorder Primer sheet header columns: name, sequence, id, type -> sequence, id, type, name
$wsPrimer = $wb.Worksheets.Item("Primer")
$wsPrimer.Range("A1").Value = "sequence"
$wsPrimer.Range("B1").Value = "id"
$wsPrimer.Range("C1").Value = "type"
$wsPrimer.Range("D1").Value = "name"
